$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Диаграмма классов" picture (InlineShapes #2, wp14:anchorId 3254E117):
#    resized, re-rotated/flipped (rot=10800000 flipH=1 flipV=1), new
#    extent/effectExtent and a new wp14:editId. We surgically rewrite just
#    that inline drawing's OOXML via Range.WordOpenXML / InsertXML (there is
#    no InlineShape.Rotation in the object model), being careful to remove
#    the auto-appended empty paragraph that InsertXML leaves behind so the
#    surrounding paragraphs are not disturbed.
# ---------------------------------------------------------------------------
$classPic = $d.InlineShapes.Item(2)
$picRange = $classPic.Range
$picStart = $picRange.Start
$picXml = $picRange.WordOpenXML

$picXml = $picXml -replace 'wp14:editId="459FB8D6"', 'wp14:editId="6CE81CEE"'
$picXml = $picXml -replace '<wp:extent cx="5123331" cy="3181350"/>', '<wp:extent cx="5565913" cy="3180715"/>'
$picXml = $picXml -replace '<wp:effectExtent l="0" t="0" r="1270" b="0"/>', '<wp:effectExtent l="0" t="0" r="0" b="635"/>'
$picXml = $picXml -replace '<a:xfrm flipV="1">', '<a:xfrm rot="10800000" flipH="1" flipV="1">'
$picXml = $picXml -replace '<a:ext cx="5162804" cy="3205861"/>', '<a:ext cx="5588130" cy="3193411"/>'

# Remove the picture's single anchor character...
$picRange.Delete()
# ...then re-insert the modified drawing at the same spot via a fresh
# (non-stale) collapsed range.
$insertRange = $d.Range($picStart, $picStart)
$insertRange.InsertXML($picXml)

# InsertXML brings along a boilerplate trailing empty paragraph from the
# XML-fragment wrapper; drop it so paragraph structure matches the original.
# (Its Range.Text is just the paragraph mark, chr(13).)
$cleanupRange = $d.Range($picStart + 1, $picStart + 2)
if ($cleanupRange.Text -eq [char]13) {
    $cleanupRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) Two pictures gain an explicit <w:rPr><w:noProof/></w:rPr> on their run
#    (communications diagram screenshot, anchorId 58083AF9; state diagram
#    screenshot, anchorId 623CDB4E). Range.NoProofing maps directly onto
#    w:noProof in the run's rPr.
# ---------------------------------------------------------------------------
foreach ($shp in $d.InlineShapes) {
    $inlineXml = $shp.Range.WordOpenXML
    if (($inlineXml -like '*wp14:anchorId="58083AF9"*') -or ($inlineXml -like '*wp14:anchorId="623CDB4E"*')) {
        $shp.Range.NoProofing = $true
    }
}

$d.Saved = $false
